$d = $word.ActiveDocument

# 1) Title line: drop the trailing "Due 12/8 (Mon) Midnight" (and the 4 tabs
#    preceding it) so the first paragraph ends right after "HW 8".
$d.Content.Find.Execute(
    "HW 8" + [char]9 + [char]9 + [char]9 + [char]9 + "Due 12/8 (Mon) Midnight",
    $false, $false, $false, $false, $false, $true, 1, $false, "HW 8", 2) | Out-Null

# 2) "ssh -l  yourCSusername" -> "ssh -l yourCSusername" (collapse the double
#    space right after "-l" to a single space).
$d.Content.Find.Execute(
    "ssh -l  yourCSusername",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "ssh -l yourCSusername", 2) | Out-Null

# 3) Rewrite the hand-in paragraph.
$d.Content.Find.Execute(
    "There is no explicit hand-in. The bomb will notify your instructor automatically after you have successfully defused it.  ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "To make sure you get credit, please submit your username and bomb number to blackboard. There is no other hand-in for this assignment. The bomb will automatically report your score as you successfully defuse it.  ",
    2) | Out-Null
